# The deck's slide master carried the "Integral" (Red Violet) theme
# colours. The authored change recolours it to the stock default
# "Office" colour palette (the palette that, before the edit, only
# lived in the otherwise-unused second theme part used by the notes
# master). Recolour the live theme - the one backing the slide master
# and therefore every slide - index by index via ThemeColorScheme,
# PowerPoint's OM entry point onto the theme part's <a:clrScheme>.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# ThemeColorScheme index -> element -> target "Office" RGB value.
# (No RGB() helper in this host, so the R+G*256+B*65536 longs are
# spelled out explicitly alongside the hex they encode.)
$colors.Colors(1).RGB  = 0x000000   # dk1       000000
$colors.Colors(2).RGB  = 0xFFFFFF   # lt1       FFFFFF
$colors.Colors(3).RGB  = 0x6A5444   # dk2       44546A
$colors.Colors(4).RGB  = 0xE6E6E7   # lt2       E7E6E6
$colors.Colors(5).RGB  = 0xD59B5B   # accent1   5B9BD5
$colors.Colors(6).RGB  = 0x317DED   # accent2   ED7D31
$colors.Colors(7).RGB  = 0xA5A5A5   # accent3   A5A5A5
$colors.Colors(8).RGB  = 0x00C0FF   # accent4   FFC000
$colors.Colors(9).RGB  = 0xC47244   # accent5   4472C4
$colors.Colors(10).RGB = 0x47AD70   # accent6   70AD47
$colors.Colors(11).RGB = 0xC16305   # hlink     0563C1
$colors.Colors(12).RGB = 0x724F95   # folHlink  954F72
